$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "QVM Cost" -> sheet1.xml
$ws2 = $wb.Worksheets.Item(2)   # "QVM All"  -> sheet2.xml

# Comment / header block appended below the existing data table on both sheets.
# Rows are populated in the same entry order as the source workbook so that
# the shared-string table indices line up: row14, row15, row13, row16, row17.
$entries = @(
    @{ Row = 14; Text = "// (Copyright) Author: anonymousresearxer" },
    @{ Row = 15; Text = "// Date created: 17/03/2022" },
    @{ Row = 13; Text = "// File name: DATASET.xlsx" },
    @{ Row = 16; Text = "// Version: Excel 2019 Pro" },
    @{ Row = 17; Text = "// Description: Scalability Measurements of submission FSE'22 1179" }
)

foreach ($entry in $entries) {
    $ws1.Cells.Item($entry.Row, 1).Value = $entry.Text
    $ws2.Cells.Item($entry.Row, 1).Value = $entry.Text
}

# Restore/update the view state: sheet1 selection on the new block,
# then sheet2 selection + activation (making "QVM All" the active tab).
$ws1.Range("A13:A17").Select()
$ws2.Range("F16").Select()
